$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 83335460
$ws.Range("I127").Value = 499
$ws.Range("J127").Value = 111113780
$ws.Range("K127").Value = 1497
$ws.Range("L127").Value = 333341340
$ws.Range("M127").Value = 3463
$ws.Range("N127").Value = -333351260
$ws.Range("H129").Value = 883.9167
$ws.Range("I129").Value = 298
$ws.Range("J129").Value = 922.9778
$ws.Range("K129").Value = 894
$ws.Range("L129").Value = 2768.9334
$ws.Range("M129").Value = 4106
$ws.Range("N129").Value = -12768.9334
$ws.Range("H133").Value = 54114.285
$ws.Range("J133").Value = 54114.285
$ws.Range("L133").Value = 54114.285
$ws.Range("N133").Value = -64234.285
$ws.Range("H137").Value = 567439.6
$ws.Range("I137").Value = 2796.818
$ws.Range("J137").Value = 903173.2
$ws.Range("K137").Value = 8390.454000000002
$ws.Range("L137").Value = 2709519.6
$ws.Range("M137").Value = -5840.454000000002
$ws.Range("N137").Value = -2714619.6
$ws.Range("H138").Value = 3465.145
$ws.Range("I138").Value = 1810.0625
$ws.Range("J138").Value = 3964.7925
$ws.Range("K138").Value = 5430.1875
$ws.Range("L138").Value = 11894.3775
$ws.Range("M138").Value = -290.1875
$ws.Range("N138").Value = -22174.3775

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18873.787
$ws.Range("I32").Value = 19933.268
$ws.Range("J32").Value = 7007.6
$ws.Range("K32").Value = 19933.268
$ws.Range("L32").Value = 7007.6
$ws.Range("M32").Value = -19646.268
$ws.Range("N32").Value = -7581.6
$ws.Range("H45").Value = 2152.6086
$ws.Range("I45").Value = 2071.9048
$ws.Range("K45").Value = 2071.9048
$ws.Range("M45").Value = -1694.9048
$ws.Range("H61").Value = 8825.583000000001
$ws.Range("I61").Value = 4890.8335
$ws.Range("J61").Value = 16695.084
$ws.Range("K61").Value = 4890.8335
$ws.Range("L61").Value = 16695.084
$ws.Range("M61").Value = -4678.8335
$ws.Range("N61").Value = -17119.084
$ws.Range("H110").Value = 2000.7858
$ws.Range("I110").Value = 1992.5834
$ws.Range("J110").Value = 2050
$ws.Range("K110").Value = 1992.5834
$ws.Range("L110").Value = 2050
$ws.Range("M110").Value = 52.41660000000002
$ws.Range("N110").Value = -6140
$ws.Range("H122").Value = 12501760
$ws.Range("I122").Value = 1950
$ws.Range("K122").Value = 5850
$ws.Range("M122").Value = -3400
$ws.Range("H132").Value = 3545.7917
$ws.Range("I132").Value = 3256.3125
$ws.Range("J132").Value = 4124.75
$ws.Range("K132").Value = 9768.9375
$ws.Range("L132").Value = 12374.25
$ws.Range("M132").Value = -7238.9375
$ws.Range("N132").Value = -17434.25
$ws.Range("H136").Value = 8825.583000000001
$ws.Range("I136").Value = 4890.8335
$ws.Range("J136").Value = 16695.084
$ws.Range("K136").Value = 14672.5005
$ws.Range("L136").Value = 50085.25199999999
$ws.Range("M136").Value = -12122.5005
$ws.Range("N136").Value = -55185.25199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2854.2
$ws.Range("I107").Value = 2636.3635
$ws.Range("J107").Value = 3453.25
$ws.Range("K107").Value = 2636.3635
$ws.Range("L107").Value = 3453.25
$ws.Range("M107").Value = -716.3634999999999
$ws.Range("N107").Value = -7293.25
$ws.Range("H114").Value = 79800
$ws.Range("J114").Value = 79800
$ws.Range("L114").Value = 79800
$ws.Range("N114").Value = -88478

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 5000
$ws.Range("J13").Value = 5000
$ws.Range("L13").Value = 5000
$ws.Range("N13").Value = -5278
$ws.Range("H31").Value = 732039.1
$ws.Range("I31").Value = 13994.467
$ws.Range("J31").Value = 1031224.44
$ws.Range("K31").Value = 13994.467
$ws.Range("L31").Value = 1031224.44
$ws.Range("M31").Value = -13699.467
$ws.Range("N31").Value = -1031814.44
$ws.Range("H34").Value = 732039.1
$ws.Range("I34").Value = 13994.467
$ws.Range("J34").Value = 1031224.44
$ws.Range("K34").Value = 13994.467
$ws.Range("L34").Value = 1031224.44
$ws.Range("M34").Value = -13792.467
$ws.Range("N34").Value = -1031628.44
$ws.Range("H43").Value = 21000
$ws.Range("J43").Value = 21000
$ws.Range("L43").Value = 21000
$ws.Range("N43").Value = -21368
$ws.Range("H58").Value = 1937868.9
$ws.Range("I58").Value = 3135981.2
$ws.Range("J58").Value = 7576.6113
$ws.Range("K58").Value = 3135981.2
$ws.Range("L58").Value = 7576.6113
$ws.Range("M58").Value = -3135778.2
$ws.Range("N58").Value = -7982.6113
$ws.Range("H101").Value = 21000
$ws.Range("J101").Value = 21000
$ws.Range("L101").Value = 21000
$ws.Range("N101").Value = -27490
$ws.Range("H102").Value = 55000
$ws.Range("J102").Value = 55000
$ws.Range("L102").Value = 55000
$ws.Range("N102").Value = -59868
$ws.Range("H132").Value = 2738.0698
$ws.Range("I132").Value = 2450.3235
$ws.Range("J132").Value = 3825.111
$ws.Range("K132").Value = 7350.970499999999
$ws.Range("L132").Value = 11475.333
$ws.Range("M132").Value = -4820.970499999999
$ws.Range("N132").Value = -16535.333
$ws.Range("H134").Value = 2663.742
$ws.Range("I134").Value = 2097.524
$ws.Range("K134").Value = 6292.572
$ws.Range("M134").Value = -3757.572
$ws.Range("H136").Value = 1937868.9
$ws.Range("I136").Value = 3135981.2
$ws.Range("J136").Value = 7576.6113
$ws.Range("K136").Value = 9407943.600000001
$ws.Range("L136").Value = 22729.8339
$ws.Range("M136").Value = -9405393.600000001
$ws.Range("N136").Value = -27829.8339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 18530064
$ws.Range("I5").Value = 483.66666
$ws.Range("J5").Value = 55589224
$ws.Range("K5").Value = 1450.99998
$ws.Range("L5").Value = 166767672
$ws.Range("M5").Value = -1338.99998
$ws.Range("N5").Value = -166767896
$ws.Range("H69").Value = 125003780
$ws.Range("J69").Value = 166671330
$ws.Range("L69").Value = 500013990
$ws.Range("N69").Value = -500015612
$ws.Range("H72").Value = 125003780
$ws.Range("J72").Value = 166671330
$ws.Range("L72").Value = 1500041970
$ws.Range("N72").Value = -1500050082
$ws.Range("H76").Value = 6997.5
$ws.Range("J76").Value = 6997.5
$ws.Range("L76").Value = 20992.5
$ws.Range("N76").Value = -21758.5
$ws.Range("H79").Value = 6997.5
$ws.Range("J79").Value = 6997.5
$ws.Range("L79").Value = 20992.5
$ws.Range("N79").Value = -23644.5
$ws.Range("H107").Value = 601.16
$ws.Range("I107").Value = 452.0625
$ws.Range("J107").Value = 738.78845
$ws.Range("K107").Value = 1356.1875
$ws.Range("L107").Value = 2216.36535
$ws.Range("M107").Value = 563.8125
$ws.Range("N107").Value = -6056.36535
$ws.Range("H131").Value = 1246.711
$ws.Range("I131").Value = 1827
$ws.Range("J131").Value = 1080.9143
$ws.Range("K131").Value = 5481
$ws.Range("L131").Value = 3242.7429
$ws.Range("M131").Value = -441
$ws.Range("N131").Value = -13322.7429
$ws.Range("H132").Value = 2790.5386
$ws.Range("I132").Value = 2030.7778
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 18277.0002
$ws.Range("L132").Value = 40500
$ws.Range("M132").Value = -15747.0002
$ws.Range("N132").Value = -45560
$ws.Range("H135").Value = 18530064
$ws.Range("I135").Value = 483.66666
$ws.Range("J135").Value = 55589224
$ws.Range("K135").Value = 4352.99994
$ws.Range("L135").Value = 500303016
$ws.Range("M135").Value = -1817.99994
$ws.Range("N135").Value = -500308086
$ws.Range("H140").Value = 1905.4546
$ws.Range("I140").Value = 1320.7142
$ws.Range("J140").Value = 5180
$ws.Range("K140").Value = 3962.1426
$ws.Range("L140").Value = 15540
$ws.Range("M140").Value = 1217.8574
$ws.Range("N140").Value = -25900
$ws.Range("H141").Value = 2403.0625
$ws.Range("J141").Value = 5979.5
$ws.Range("L141").Value = 17938.5
$ws.Range("N141").Value = -28298.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 13979.333
$ws.Range("I122").Value = 26201.5
$ws.Range("J122").Value = 4201.6
$ws.Range("K122").Value = 78604.5
$ws.Range("L122").Value = 12604.8
$ws.Range("M122").Value = -76154.5
$ws.Range("N122").Value = -17504.8
$ws.Range("H129").Value = 49992.668
$ws.Range("J129").Value = 49992.668
$ws.Range("L129").Value = 49992.668
$ws.Range("N129").Value = -59992.668
$ws.Range("H132").Value = 14766.941
$ws.Range("I132").Value = 16917.715
$ws.Range("J132").Value = 13261.4
$ws.Range("K132").Value = 50753.145
$ws.Range("L132").Value = 39784.2
$ws.Range("M132").Value = -48223.145
$ws.Range("N132").Value = -44844.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3487.125
$ws.Range("I7").Value = 3656.8572
$ws.Range("J7").Value = 2299
$ws.Range("K7").Value = 3656.8572
$ws.Range("L7").Value = 2299
$ws.Range("M7").Value = -3544.8572
$ws.Range("N7").Value = -2523
$ws.Range("H22").Value = 400.07693
$ws.Range("I22").Value = 399.44446
$ws.Range("J22").Value = 401.5
$ws.Range("K22").Value = 399.44446
$ws.Range("L22").Value = 401.5
$ws.Range("M22").Value = -104.44446
$ws.Range("N22").Value = -991.5
$ws.Range("H27").Value = 400.07693
$ws.Range("I27").Value = 399.44446
$ws.Range("J27").Value = 401.5
$ws.Range("K27").Value = 399.44446
$ws.Range("L27").Value = 401.5
$ws.Range("M27").Value = -292.44446
$ws.Range("N27").Value = -615.5
$ws.Range("H40").Value = 3284.16
$ws.Range("I40").Value = 3243.75
$ws.Range("J40").Value = 3356
$ws.Range("K40").Value = 3243.75
$ws.Range("L40").Value = 3356
$ws.Range("M40").Value = -3107.75
$ws.Range("N40").Value = -3628
$ws.Range("H123").Value = 57300
$ws.Range("J123").Value = 57300
$ws.Range("L123").Value = 57300
$ws.Range("N123").Value = -67100
$ws.Range("H126").Value = 3487.125
$ws.Range("I126").Value = 3656.8572
$ws.Range("J126").Value = 2299
$ws.Range("K126").Value = 10970.5716
$ws.Range("L126").Value = 6897
$ws.Range("M126").Value = -8500.571599999999
$ws.Range("N126").Value = -11837

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3307.0667
$ws.Range("J107").Value = 3801.25
$ws.Range("L107").Value = 11403.75
$ws.Range("N107").Value = -15243.75
$ws.Range("H122").Value = 1750
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -16900
$ws.Range("H126").Value = 1844.6666
$ws.Range("I126").Value = 1813.6
$ws.Range("K126").Value = 5440.799999999999
$ws.Range("M126").Value = -2970.799999999999
$ws.Range("H132").Value = 2127.2727
$ws.Range("I132").Value = 1992.5714
$ws.Range("J132").Value = 2881.6
$ws.Range("K132").Value = 5977.7142
$ws.Range("L132").Value = 8644.799999999999
$ws.Range("M132").Value = -3447.7142
$ws.Range("N132").Value = -13704.8
